$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 461, shifting existing rows 461:545 down to 462:546
$ws.Rows.Item(461).Insert()

# Populate the newly inserted row 461 with the new record's data
$ws.Cells.Item(461, 1).Value = 3
$ws.Cells.Item(461, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(461, 3).Value = "Coquimbo"
$ws.Cells.Item(461, 4).Value = 44995
$ws.Cells.Item(461, 5).Value = 5
$ws.Cells.Item(461, 6).Value = 100112031
$ws.Cells.Item(461, 7).Value = "Poroto verde"
$ws.Cells.Item(461, 8).Value = "Magnum"
$ws.Cells.Item(461, 9).Value = "Primera"
$ws.Cells.Item(461, 10).Value = 73
$ws.Cells.Item(461, 11).Value = 25000
$ws.Cells.Item(461, 12).Value = 26000
$ws.Cells.Item(461, 13).Value = 25521
$ws.Cells.Item(461, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(461, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(461, 16).Value = 1021
$ws.Cells.Item(461, 17).Value = 25
$ws.Cells.Item(461, 18).Value = "Hortaliza"
